{"js": "// Replace the date and every division-problem answer cell with the new values.\n// Each old value is unique in the document, so a targeted search+replace per\n// value is safe and precise.\nconst replacements = [\n  [\"2023-08-23 Wednesday\", \"2023-08-24 Thursday\"],\n  [\"29\u00f74=7, 1\", \"19\u00f72=9, 1\"],\n  [\"25\u00f75=5, 0\", \"43\u00f77=6, 1\"],\n  [\"50\u00f77=7, 1\", \"13\u00f77=1, 6\"],\n  [\"33\u00f76=5, 3\", \"66\u00f73=22, 0\"],\n  [\"17\u00f78=2, 1\", \"34\u00f72=17, 0\"],\n  [\"81\u00f73=27, 0\", \"87\u00f78=10, 7\"],\n  [\"45\u00f78=5, 5\", \"55\u00f73=18, 1\"],\n  [\"48\u00f77=6, 6\", \"10\u00f74=2, 2\"],\n  [\"48\u00f75=9, 3\", \"55\u00f72=27, 1\"],\n  [\"43\u00f78=5, 3\", \"66\u00f79=7, 3\"],\n  [\"79\u00f75=15, 4\", \"10\u00f77=1, 3\"],\n  [\"30\u00f76=5, 0\", \"49\u00f77=7, 0\"],\n  [\"49\u00f76=8, 1\", \"35\u00f76=5, 5\"],\n  [\"63\u00f78=7, 7\", \"95\u00f72=47, 1\"],\n  [\"86\u00f76=14, 2\", \"70\u00f78=8, 6\"],\n  [\"93\u00f76=15, 3\", \"23\u00f72=11, 1\"],\n  [\"60\u00f78=7, 4\", \"71\u00f74=17, 3\"],\n  [\"48\u00f76=8, 0\", \"33\u00f75=6, 3\"],\n  [\"16\u00f78=2, 0\", \"96\u00f76=16, 0\"],\n  [\"41\u00f74=10, 1\", \"25\u00f73=8, 1\"],\n  [\"49\u00f78=6, 1\", \"85\u00f75=17, 0\"],\n  [\"51\u00f75=10, 1\", \"24\u00f72=12, 0\"],\n  [\"59\u00f74=14, 3\", \"49\u00f73=16, 1\"],\n  [\"27\u00f76=4, 3\", \"42\u00f77=6, 0\"],\n  [\"76\u00f76=12, 4\", \"12\u00f73=4, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date and every division-problem answer cell with the new\n# values. Each old value is unique in the document, so a Find/Replace pass\n# per value precisely targets the right run.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2023-08-23 Wednesday\", \"2023-08-24 Thursday\"),\n    @(\"29\u00f74=7, 1\", \"19\u00f72=9, 1\"),\n    @(\"25\u00f75=5, 0\", \"43\u00f77=6, 1\"),\n    @(\"50\u00f77=7, 1\", \"13\u00f77=1, 6\"),\n    @(\"33\u00f76=5, 3\", \"66\u00f73=22, 0\"),\n    @(\"17\u00f78=2, 1\", \"34\u00f72=17, 0\"),\n    @(\"81\u00f73=27, 0\", \"87\u00f78=10, 7\"),\n    @(\"45\u00f78=5, 5\", \"55\u00f73=18, 1\"),\n    @(\"48\u00f77=6, 6\", \"10\u00f74=2, 2\"),\n    @(\"48\u00f75=9, 3\", \"55\u00f72=27, 1\"),\n    @(\"43\u00f78=5, 3\", \"66\u00f79=7, 3\"),\n    @(\"79\u00f75=15, 4\", \"10\u00f77=1, 3\"),\n    @(\"30\u00f76=5, 0\", \"49\u00f77=7, 0\"),\n    @(\"49\u00f76=8, 1\", \"35\u00f76=5, 5\"),\n    @(\"63\u00f78=7, 7\", \"95\u00f72=47, 1\"),\n    @(\"86\u00f76=14, 2\", \"70\u00f78=8, 6\"),\n    @(\"93\u00f76=15, 3\", \"23\u00f72=11, 1\"),\n    @(\"60\u00f78=7, 4\", \"71\u00f74=17, 3\"),\n    @(\"48\u00f76=8, 0\", \"33\u00f75=6, 3\"),\n    @(\"16\u00f78=2, 0\", \"96\u00f76=16, 0\"),\n    @(\"41\u00f74=10, 1\", \"25\u00f73=8, 1\"),\n    @(\"49\u00f78=6, 1\", \"85\u00f75=17, 0\"),\n    @(\"51\u00f75=10, 1\", \"24\u00f72=12, 0\"),\n    @(\"59\u00f74=14, 3\", \"49\u00f73=16, 1\"),\n    @(\"27\u00f76=4, 3\", \"42\u00f77=6, 0\"),\n    @(\"76\u00f76=12, 4\", \"12\u00f73=4, 0\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
